$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H18").Value = 4624.25
$ws.Range("J18").Value = 6000
$ws.Range("L18").Value = 6000
$ws.Range("N18").Value = -6568
$ws.Range("H19").Value = 1869.8572
$ws.Range("J19").Value = 2997.3333
$ws.Range("L19").Value = 2997.3333
$ws.Range("N19").Value = -3347.3333
$ws.Range("H28").Value = 828.5
$ws.Range("I28").Value = 828.5
$ws.Range("K28").Value = 828.5
$ws.Range("M28").Value = -343.5
$ws.Range("H40").Value = 4040.8696
$ws.Range("I40").Value = 2535.4285
$ws.Range("K40").Value = 2535.4285
$ws.Range("M40").Value = -2360.4285
$ws.Range("H51").Value = 3998
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 4497.5
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 4497.5
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -5465.5
$ws.Range("H61").Value = 1299.5
$ws.Range("I61").Value = 1299.5
$ws.Range("K61").Value = 3898.5
$ws.Range("M61").Value = -3726.5
$ws.Range("H99").Value = 1721.3334
$ws.Range("I99").Value = 1264
$ws.Range("K99").Value = 3792
$ws.Range("M99").Value = -2294
$ws.Range("H112").Value = 1811.5454
$ws.Range("I112").Value = 1399.5
$ws.Range("J112").Value = 1903.1111
$ws.Range("K112").Value = 4198.5
$ws.Range("L112").Value = 5709.3333
$ws.Range("M112").Value = -3090.5
$ws.Range("N112").Value = -7925.3333
$ws.Range("H116").Value = 9644.556
$ws.Range("I116").Value = 9749.75
$ws.Range("J116").Value = 9560.4
$ws.Range("K116").Value = 9749.75
$ws.Range("L116").Value = 9560.4
$ws.Range("M116").Value = -6307.75
$ws.Range("N116").Value = -16444.4
$ws.Range("H135").Value = 4310.2856
$ws.Range("I135").Value = 3918
$ws.Range("K135").Value = 35262
$ws.Range("M135").Value = -32727

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 6890.4546
$ws.Range("I61").Value = 6858.4116
$ws.Range("J61").Value = 6999.4
$ws.Range("K61").Value = 6858.4116
$ws.Range("L61").Value = 6999.4
$ws.Range("M61").Value = -6646.4116
$ws.Range("N61").Value = -7423.4
$ws.Range("H104").Value = 27737.166
$ws.Range("J104").Value = 26243.8
$ws.Range("L104").Value = 26243.8
$ws.Range("N104").Value = -33231.8
$ws.Range("H110").Value = 2762.0625
$ws.Range("J110").Value = 1013
$ws.Range("L110").Value = 1013
$ws.Range("N110").Value = -5103
$ws.Range("H136").Value = 6890.4546
$ws.Range("I136").Value = 6858.4116
$ws.Range("J136").Value = 6999.4
$ws.Range("K136").Value = 20575.2348
$ws.Range("L136").Value = 20998.2
$ws.Range("M136").Value = -18025.2348
$ws.Range("N136").Value = -26098.2

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H2").Value = 51518.6
$ws.Range("J2").Value = 50648.25
$ws.Range("L2").Value = 50648.25
$ws.Range("N2").Value = -50874.25
$ws.Range("H107").Value = 2519.15
$ws.Range("I107").Value = 1180.3334
$ws.Range("J107").Value = 4527.375
$ws.Range("K107").Value = 1180.3334
$ws.Range("L107").Value = 4527.375
$ws.Range("M107").Value = 739.6666
$ws.Range("N107").Value = -8367.375

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 2166.6667
$ws.Range("I16").Value = 1750
$ws.Range("K16").Value = 1750
$ws.Range("M16").Value = -1463
$ws.Range("H99").Value = 2800
$ws.Range("I99").Value = 2800
$ws.Range("K99").Value = 2800
$ws.Range("M99").Value = -1302
$ws.Range("H107").Value = 994.7143
$ws.Range("I107").Value = 994.7143
$ws.Range("K107").Value = 994.7143
$ws.Range("M107").Value = 925.2857
$ws.Range("H113").Value = 2166.6667
$ws.Range("I113").Value = 1750
$ws.Range("K113").Value = 1750
$ws.Range("M113").Value = 420
$ws.Range("H126").Value = 2800
$ws.Range("I126").Value = 2800
$ws.Range("K126").Value = 8400
$ws.Range("M126").Value = -5930

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H101").Value = 10056.909
$ws.Range("J101").Value = 10056.909
$ws.Range("L101").Value = 30170.727
$ws.Range("N101").Value = -35038.727
$ws.Range("H137").Value = 12926.95
$ws.Range("J137").Value = 13502.353
$ws.Range("L137").Value = 40507.05899999999
$ws.Range("N137").Value = -50707.05899999999
$ws.Range("H140").Value = 3008.9412
$ws.Range("I140").Value = 2166.077
$ws.Range("K140").Value = 6498.231000000001
$ws.Range("M140").Value = -1318.231000000001
$ws.Range("H141").Value = 111115780
$ws.Range("J141").Value = 6933
$ws.Range("L141").Value = 20799
$ws.Range("N141").Value = -31159

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H41").Value = 3999.6667
$ws.Range("I41").Value = 3499.5
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 3499.5
$ws.Range("L41").Value = 5000
$ws.Range("M41").Value = -3144.5
$ws.Range("N41").Value = -5710
$ws.Range("H102").Value = 6465.5
$ws.Range("I102").Value = 6389.5713
$ws.Range("K102").Value = 6389.5713
$ws.Range("M102").Value = -4767.5713
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988
$ws.Range("H113").Value = 13151.895
$ws.Range("I113").Value = 19898.666
$ws.Range("J113").Value = 10038
$ws.Range("K113").Value = 19898.666
$ws.Range("L113").Value = 10038
$ws.Range("M113").Value = -17728.666
$ws.Range("N113").Value = -14378
$ws.Range("H122").Value = 11857.417
$ws.Range("J122").Value = 5399.75
$ws.Range("L122").Value = 16199.25
$ws.Range("N122").Value = -21099.25

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 4107.7144
$ws.Range("I7").Value = 4250.6665
$ws.Range("J7").Value = 3250
$ws.Range("K7").Value = 4250.6665
$ws.Range("L7").Value = 3250
$ws.Range("M7").Value = -4138.6665
$ws.Range("N7").Value = -3474
$ws.Range("H40").Value = 6228.5835
$ws.Range("I40").Value = 4497.6665
$ws.Range("K40").Value = 4497.6665
$ws.Range("M40").Value = -4361.6665
$ws.Range("H68").Value = 7999.6665
$ws.Range("I68").Value = 5499.25
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 5499.25
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -4750.25
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 7999.6665
$ws.Range("I71").Value = 5499.25
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 27496.25
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -23752.25
$ws.Range("N71").Value = -57488
$ws.Range("H100").Value = 6579.577
$ws.Range("I100").Value = 5959.3887
$ws.Range("K100").Value = 5959.3887
$ws.Range("M100").Value = -5418.3887
$ws.Range("H126").Value = 4107.7144
$ws.Range("I126").Value = 4250.6665
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 12751.9995
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -10281.9995
$ws.Range("N126").Value = -14690
$ws.Range("H136").Value = 24685.625
$ws.Range("I136").Value = 4507.769
$ws.Range("J136").Value = 48532.184
$ws.Range("K136").Value = 13523.307
$ws.Range("L136").Value = 145596.552
$ws.Range("M136").Value = -10973.307
$ws.Range("N136").Value = -150696.552

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = 2109
$ws.Range("J81").Value = 3143.25
$ws.Range("L81").Value = 6286.5
$ws.Range("N81").Value = -8408.5
$ws.Range("H84").Value = 2109
$ws.Range("J84").Value = 3143.25
$ws.Range("L84").Value = 31432.5
$ws.Range("N84").Value = -42040.5
$ws.Range("H107").Value = 621.2778
$ws.Range("I107").Value = 525.38464
$ws.Range("K107").Value = 1576.15392
$ws.Range("M107").Value = 343.84608
$ws.Range("H113").Value = 576.64514
$ws.Range("I113").Value = 460.61905
$ws.Range("K113").Value = 1381.85715
$ws.Range("M113").Value = 788.14285
$ws.Range("H122").Value = 3488.6047
$ws.Range("I122").Value = 2157.6428
$ws.Range("J122").Value = 4131.1377
$ws.Range("K122").Value = 6472.928400000001
$ws.Range("L122").Value = 12393.4131
$ws.Range("M122").Value = -4022.928400000001
$ws.Range("N122").Value = -17293.4131
$ws.Range("H126").Value = 3151.85
$ws.Range("I126").Value = 2499.1667
$ws.Range("J126").Value = 4130.875
$ws.Range("K126").Value = 7497.500100000001
$ws.Range("L126").Value = 12392.625
$ws.Range("M126").Value = -5027.500100000001
$ws.Range("N126").Value = -17332.625
$ws.Range("H132").Value = 3021.1924
$ws.Range("I132").Value = 2527.05
$ws.Range("J132").Value = 4668.3335
$ws.Range("K132").Value = 7581.150000000001
$ws.Range("L132").Value = 14005.0005
$ws.Range("M132").Value = -5051.150000000001
$ws.Range("N132").Value = -19065.0005
